# CameraChangePos at the beginning
# Insert a new "cameraPos" field/column before the existing "nextCampId"
# column (old column I), shifting it (and the data below it) one column
# to the right (to column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; this shifts the old I (nextCampId) column to J
# and inherits formatting from the column to the left, same as Excel's
# normal "Insert Sheet Columns" behavior.
$ws.Columns("I:I").Insert()

# Row 1 holds the "type" marker for each field; the new cameraPos field is
# a string, like the other position/name fields.
$ws.Range("I1").Value = "string"

# Row 2 holds the field name.
$ws.Range("I2").Value = "cameraPos"

# Rows 3 and 4 hold the data values for the two map rows.
$ws.Range("I3").Value = "[-10,30]"
$ws.Range("I4").Value = "[-10,30]"

# Match column width of the new column to its neighbors (E:H).
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth

# Update the view state (zoom + selection) to match what was left active.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("I4").Select()
